$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 516
$ws.Range("I2").Value = 1320
$ws.Range("K2").Value = 41
$ws.Range("L2").Value = 1405
$ws.Range("M2").Value = 66
$ws.Range("N2").Value = 943
$ws.Range("P2").Value = 16
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 67
$ws.Range("S2").Value = 589
$ws.Range("T2").Value = 964
$ws.Range("U2").Value = 65
$ws.Range("V2").Value = 8298
$ws.Range("X2").Value = 8306
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 125
$ws.Range("AA2").Value = 65
